# Scheduled runner update: refresh cached market-board profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) for the
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1249
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1249
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1249
$ws.Range("N32").Value = -1901

$ws.Range("H70").Value = 3283.3333
$ws.Range("I70").Value = 2400
$ws.Range("J70").Value = 4166.6665
$ws.Range("K70").Value = 7200
$ws.Range("L70").Value = 12499.9995
$ws.Range("M70").Value = -6930
$ws.Range("N70").Value = -13039.9995

$ws.Range("H73").Value = 3283.3333
$ws.Range("I73").Value = 2400
$ws.Range("J73").Value = 4166.6665
$ws.Range("K73").Value = 7200
$ws.Range("L73").Value = 12499.9995
$ws.Range("M73").Value = -6264
$ws.Range("N73").Value = -14371.9995

$ws.Range("H76").Value = 3219.8
$ws.Range("I76").Value = 3024.75
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3024.75
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2709.75
$ws.Range("N76").Value = -4630

$ws.Range("H79").Value = 3219.8
$ws.Range("I79").Value = 3024.75
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3024.75
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -1932.75
$ws.Range("N79").Value = -6184

$ws.Range("H99").Value = 66667920
$ws.Range("I99").Value = 111111200
$ws.Range("J99").Value = 2999.5
$ws.Range("K99").Value = 333333600
$ws.Range("L99").Value = 8998.5
$ws.Range("M99").Value = -333332102
$ws.Range("N99").Value = -11994.5

$ws.Range("H116").Value = 7502.5
$ws.Range("I116").Value = 7502.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 7502.5
$ws.Range("L116").ClearContents()
$ws.Range("M116").Value = -4060.5
$ws.Range("N116").ClearContents()

$ws.Range("H137").Value = 1267.3334
$ws.Range("I137").Value = 921
$ws.Range("J137").Value = 2999
$ws.Range("K137").Value = 2763
$ws.Range("L137").Value = 8997
$ws.Range("M137").Value = -213
$ws.Range("N137").Value = -14097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 100
$ws.Range("N2").Value = -326

$ws.Range("H32").Value = 969.1667
$ws.Range("I32").Value = 969.1667
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 969.1667
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -682.1667

$ws.Range("H76").Value = 28665
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 28665
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 28665
$ws.Range("N76").Value = -29341

$ws.Range("H79").Value = 28665
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 28665
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 28665
$ws.Range("N79").Value = -31005

$ws.Range("H88").Value = 4033.75
$ws.Range("I88").Value = 2878.5
$ws.Range("J88").Value = 7499.5
$ws.Range("K88").Value = 2878.5
$ws.Range("L88").Value = 7499.5
$ws.Range("M88").Value = -2472.5
$ws.Range("N88").Value = -8311.5

$ws.Range("H91").Value = 4033.75
$ws.Range("I91").Value = 2878.5
$ws.Range("J91").Value = 7499.5
$ws.Range("K91").Value = 2878.5
$ws.Range("L91").Value = 7499.5
$ws.Range("M91").Value = -1474.5
$ws.Range("N91").Value = -10307.5

$ws.Range("H102").Value = 11724834
$ws.Range("I102").Value = 918918.25
$ws.Range("J102").Value = 33336666
$ws.Range("K102").Value = 918918.25
$ws.Range("L102").Value = 33336666
$ws.Range("M102").Value = -917296.25
$ws.Range("N102").Value = -33339910

$ws.Range("H116").Value = 100
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 100
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 100
$ws.Range("N116").Value = -4688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100
$ws.Range("N3").Value = -328

$ws.Range("H22").Value = 751.6667
$ws.Range("I22").Value = 981.5714
$ws.Range("J22").Value = 429.8
$ws.Range("K22").Value = 981.5714
$ws.Range("L22").Value = 429.8
$ws.Range("M22").Value = -808.5714
$ws.Range("N22").Value = -775.8

$ws.Range("H86").Value = 3883
$ws.Range("I86").Value = 3883
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3883
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2760

$ws.Range("H89").Value = 3883
$ws.Range("I89").Value = 3883
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 19415
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -13799

$ws.Range("H99").Value = 4649.467
$ws.Range("I99").Value = 4553
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 4553
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -3055
$ws.Range("N99").Value = -8996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 62523
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 62523
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 62523
$ws.Range("N74").Value = -64271

$ws.Range("H77").Value = 62523
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 62523
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 187569
$ws.Range("N77").Value = -196305

$ws.Range("H80").Value = 50064
$ws.Range("I80").Value = 55000
$ws.Range("J80").Value = 45128
$ws.Range("K80").Value = 55000
$ws.Range("L80").Value = 45128
$ws.Range("M80").Value = -53877
$ws.Range("N80").Value = -47374

$ws.Range("H83").Value = 50064
$ws.Range("I83").Value = 55000
$ws.Range("J83").Value = 45128
$ws.Range("K83").Value = 165000
$ws.Range("L83").Value = 135384
$ws.Range("M83").Value = -159384
$ws.Range("N83").Value = -146616

$ws.Range("H86").Value = 333350000
$ws.Range("I86").Value = 500010000
$ws.Range("J86").Value = 30000
$ws.Range("K86").Value = 500010000
$ws.Range("L86").Value = 30000
$ws.Range("M86").Value = -500008877
$ws.Range("N86").Value = -32246

$ws.Range("H89").Value = 333350000
$ws.Range("I89").Value = 500010000
$ws.Range("J89").Value = 30000
$ws.Range("K89").Value = 2500050000
$ws.Range("L89").Value = 150000
$ws.Range("M89").Value = -2500044384
$ws.Range("N89").Value = -161232

$ws.Range("H103").Value = 43295
$ws.Range("I103").Value = 43295
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 43295
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -42123

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 399.25
$ws.Range("I2").Value = 431
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 2586
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -2473
$ws.Range("N2").Value = -526

$ws.Range("H33").Value = 535.375
$ws.Range("I33").Value = 617.2
$ws.Range("J33").Value = 399
$ws.Range("K33").Value = 3703.2
$ws.Range("L33").Value = 2394
$ws.Range("M33").Value = -3420.2
$ws.Range("N33").Value = -2960

$ws.Range("H100").Value = 2476.923
$ws.Range("I100").Value = 2476.923
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 7430.768999999999
$ws.Range("L100").ClearContents()
$ws.Range("M100").Value = -6619.768999999999
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1770
$ws.Range("I132").Value = 1699.5
$ws.Range("J132").Value = 1911
$ws.Range("K132").Value = 5098.5
$ws.Range("L132").Value = 5733
$ws.Range("M132").Value = -2568.5
$ws.Range("N132").Value = -10793

$ws.Range("H136").Value = 25000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 25000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 75000
$ws.Range("N136").Value = -80100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3575.8667
$ws.Range("I22").Value = 3348.6667
$ws.Range("J22").Value = 3916.6667
$ws.Range("K22").Value = 3348.6667
$ws.Range("L22").Value = 3916.6667
$ws.Range("M22").Value = -3053.6667
$ws.Range("N22").Value = -4506.6667

$ws.Range("H27").Value = 3575.8667
$ws.Range("I27").Value = 3348.6667
$ws.Range("J27").Value = 3916.6667
$ws.Range("K27").Value = 3348.6667
$ws.Range("L27").Value = 3916.6667
$ws.Range("M27").Value = -3241.6667
$ws.Range("N27").Value = -4130.6667

$ws.Range("H64").Value = 35489.8
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 35489.8
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 35489.8
$ws.Range("N64").Value = -35939.8

$ws.Range("H67").Value = 35489.8
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 35489.8
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 35489.8
$ws.Range("N67").Value = -37049.8

$ws.Range("H132").Value = 2463.5454
$ws.Range("I132").Value = 2285.5
$ws.Range("J132").Value = 2938.3333
$ws.Range("K132").Value = 6856.5
$ws.Range("L132").Value = 8814.999899999999
$ws.Range("M132").Value = -4326.5
$ws.Range("N132").Value = -13874.9999

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 20027.75
$ws.Range("I52").Value = 20027.75
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 20027.75
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -19801.75

$ws.Range("H68").Value = 49333.332
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 49333.332
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 49333.332
$ws.Range("N68").Value = -50955.332

$ws.Range("H71").Value = 49333.332
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 49333.332
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 147999.996
$ws.Range("N71").Value = -156111.996

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").ClearContents()
